$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.877.90'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '2.273.45'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.60'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.30'
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.73'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '2.625.39'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.39'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("D16").Value = '2.277.80'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.783'
$ws.Range("E17").Value = '  +3.44%  '
$ws.Range("D18").Value = '41.808.63'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +5.30%  '
$ws.Range("D20").Value = '0.0₃0917'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.26'
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("E25").Value = '  +1.96%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.05'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.81'
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.32'
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.38'
$ws.Range("E32").Value = '  +4.40%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.33'
$ws.Range("E36").Value = '  +4.49%  '
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.006.19'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.60'
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("E44").Value = '  +10.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0283'
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.28'
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.49'
$ws.Range("E48").Value = '  +3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '72.74'
$ws.Range("E49").Value = '  +3.03%  '
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("E51").Value = '  +0.22%  '
